$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.133750915527344
$ws.Range("B1").Value = 6.350811958312988
$ws.Range("C1").Value = 6.063872814178467
$ws.Range("D1").Value = 6.477526187896729
$ws.Range("E1").Value = 5.340754508972168
